# The four observation rows (4-7) on the "Artfynd" sheet get their
# identifying/species data cyclically rotated:
#   row4's data moves to row7
#   row5's data moves to row4
#   row6's data moves to row5
#   row7's data moves to row6
# (Location/date/observer columns are identical across these rows so they
# do not need to be touched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually differ between the four rows and therefore need
# to be rewritten.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Snapshot current values for rows 4-7 before any writes occur.
$orig = @{}
foreach ($r in 4..7) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range($col + $r).Value2
    }
    $orig[$r]["M"] = $ws.Range("M" + $r).Value2
}

# Map: new row -> source row supplying its data.
$source = @{ 4 = 5; 5 = 6; 6 = 7; 7 = 4 }

foreach ($r in 4..7) {
    $src = $source[$r]
    foreach ($col in $cols) {
        $ws.Range($col + $r).Value2 = $orig[$src][$col]
    }

    $mVal = $orig[$src]["M"]
    if ($mVal -eq $null -or $mVal -eq "") {
        $ws.Range("M" + $r).ClearContents()
    } else {
        $ws.Range("M" + $r).Value2 = $mVal
    }
}
